{"js": "// Update the multiplication problems in the practice table.\n// Each entry is [oldText, newText]; applied in document order so that\n// the one value that is both an old and a new text (873\u00d75=) resolves\n// unambiguously (the earlier 873\u00d75= is replaced with 883\u00d75= before the\n// later 955\u00d72= is turned into 873\u00d75=).\nconst replacements = [\n  [\"945\u00d74=\", \"794\u00d75=\"],\n  [\"541\u00d74=\", \"419\u00d77=\"],\n  [\"263\u00d78=\", \"634\u00d73=\"],\n  [\"958\u00d76=\", \"426\u00d76=\"],\n  [\"251\u00d79=\", \"141\u00d77=\"],\n  [\"543\u00d78=\", \"596\u00d78=\"],\n  [\"279\u00d73=\", \"842\u00d73=\"],\n  [\"217\u00d76=\", \"763\u00d74=\"],\n  [\"498\u00d73=\", \"635\u00d76=\"],\n  [\"693\u00d73=\", \"661\u00d77=\"],\n  [\"249\u00d77=\", \"569\u00d73=\"],\n  [\"426\u00d75=\", \"591\u00d77=\"],\n  [\"990\u00d76=\", \"748\u00d76=\"],\n  [\"873\u00d75=\", \"883\u00d75=\"],\n  [\"866\u00d74=\", \"422\u00d74=\"],\n  [\"577\u00d76=\", \"749\u00d77=\"],\n  [\"746\u00d78=\", \"614\u00d74=\"],\n  [\"249\u00d72=\", \"928\u00d73=\"],\n  [\"601\u00d73=\", \"720\u00d78=\"],\n  [\"570\u00d76=\", \"508\u00d76=\"],\n  [\"169\u00d76=\", \"787\u00d72=\"],\n  [\"287\u00d77=\", \"329\u00d73=\"],\n  [\"991\u00d76=\", \"809\u00d72=\"],\n  [\"778\u00d72=\", \"151\u00d74=\"],\n  [\"955\u00d72=\", \"873\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace only the first occurrence (each search term is unique at the\n  // time it is searched for because we work top-to-bottom through the\n  // document in the same order the originals appear).\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the multiplication problems in the practice table.\n# Each pair is (oldText, newText), listed in the order the original\n# values appear in the document. wdReplaceOne (1) replaces only the\n# first match on each call, so running the pairs strictly in document\n# order keeps things unambiguous even though one new value (873\u00d75=)\n# equals an earlier old value: by the time we process\n# \"955\u00d72=\" -> \"873\u00d75=\", the earlier cell's \"873\u00d75=\" has already been\n# turned into \"883\u00d75=\", so there is nothing stray left for Find to hit.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"945\u00d74=\", \"794\u00d75=\"),\n    @(\"541\u00d74=\", \"419\u00d77=\"),\n    @(\"263\u00d78=\", \"634\u00d73=\"),\n    @(\"958\u00d76=\", \"426\u00d76=\"),\n    @(\"251\u00d79=\", \"141\u00d77=\"),\n    @(\"543\u00d78=\", \"596\u00d78=\"),\n    @(\"279\u00d73=\", \"842\u00d73=\"),\n    @(\"217\u00d76=\", \"763\u00d74=\"),\n    @(\"498\u00d73=\", \"635\u00d76=\"),\n    @(\"693\u00d73=\", \"661\u00d77=\"),\n    @(\"249\u00d77=\", \"569\u00d73=\"),\n    @(\"426\u00d75=\", \"591\u00d77=\"),\n    @(\"990\u00d76=\", \"748\u00d76=\"),\n    @(\"873\u00d75=\", \"883\u00d75=\"),\n    @(\"866\u00d74=\", \"422\u00d74=\"),\n    @(\"577\u00d76=\", \"749\u00d77=\"),\n    @(\"746\u00d78=\", \"614\u00d74=\"),\n    @(\"249\u00d72=\", \"928\u00d73=\"),\n    @(\"601\u00d73=\", \"720\u00d78=\"),\n    @(\"570\u00d76=\", \"508\u00d76=\"),\n    @(\"169\u00d76=\", \"787\u00d72=\"),\n    @(\"287\u00d77=\", \"329\u00d73=\"),\n    @(\"991\u00d76=\", \"809\u00d72=\"),\n    @(\"778\u00d72=\", \"151\u00d74=\"),\n    @(\"955\u00d72=\", \"873\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdReplaceOne = 1: replace only the first occurrence found, keeping\n    # the overall operation order-safe for the single old/new value overlap.\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
